# Ticket 50: Add a "varStatus" demo worksheet to ForTagTemplate.xlsx,
# showing the new varStatus attribute / RangedLoopTagStatus (start/end/step)
# for nested jt:for loops.

$wb = $excel.ActiveWorkbook

# --- Add the new sheet as the last tab -------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "varStatus"

# --- Header row (row 1): reuse the existing bold/blue/bordered header style
# by copying the format from the "Multiplication" sheet's header cells.
$headerSource = $wb.Worksheets.Item("Multiplication")
$headerSource.Range("A1:B1").Copy()
$ws.Range("A1:H1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A1").Value = "X"
$ws.Range("B1").Value = "Y"
$ws.Range("C1").Value = "startX"
$ws.Range("D1").Value = "endX"
$ws.Range("E1").Value = "stepX"
$ws.Range("F1").Value = "startY"
$ws.Range("G1").Value = "endY"
$ws.Range("H1").Value = "stepY"

# --- Body row (row 2): plain thin border around each cell, no fill/font change.
$body = $ws.Range("A2:H2")
$body.Borders.LineStyle = 1
$body.Borders.Weight = 2

$ws.Range("A2").Value = '<jt:for var="x" start="1" end="5" varStatus="xs"><jt:for var="y" start="5" end="${x}" step="-1" varStatus="ys">${x}'
$ws.Range("B2").Value = '${y}'
$ws.Range("C2").Value = '${xs.start}'
$ws.Range("D2").Value = '${xs.end}'
$ws.Range("E2").Value = '${xs.step}'
$ws.Range("F2").Value = '${ys.start}'
$ws.Range("G2").Value = '${ys.end}'
$ws.Range("H2").Value = '${ys.step}</jt:for></jt:for>'

# Match the other template sheets' portrait page setup.
$ws.PageSetup.Orientation = 1
